$d = $word.ActiveDocument

# The target paragraph (1) currently ends the document, and the COM-interop
# host always re-merges a "last paragraph in the document" mark's pPr back
# onto whatever ends up last, so InsertXML alone cannot drop <w:jc>. Work
# around that by temporarily appending an extra paragraph, rewriting
# paragraph 1 (now no longer last) via InsertXML, then deleting the helper
# paragraph again.
$full = $d.Range(0, $d.Content.End)
$full.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Вид </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}} {{</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Декан </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$r1.InsertXML($xml) | Out-Null

# Drop the helper paragraph we appended above.
$d.Paragraphs.Item(2).Range.Delete() | Out-Null
